$wb = $excel.ActiveWorkbook

# Added 2020 results + clean up:
#  - The "Assumptions" scratch sheet (just held a tiny epsilon constant) is no
#    longer needed, so remove it. EXPORT slides up to take its place in the
#    tab order and related workbook bookkeeping (sheet ids / defined-name
#    refs / shared strings) is recomputed automatically by Excel on save.
$excel.DisplayAlerts = $false
$wb.Sheets.Item("Assumptions").Delete()

# Leave the EXPORT tab selected/active (it is now the last, 4th, tab).
$wb.Sheets.Item("EXPORT").Activate()

# Lock down all of the remaining sheets (Protect Sheet, no password) so
# casual users can't disturb the formulas/layout.
foreach ($sheetName in @("Election Results by State", "Uncontested Races", "Uncontested by State PIVOT", "EXPORT")) {
    $ws = $wb.Sheets.Item($sheetName)
    $ws.Protect($null, $true, $true, $true)
}
